$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.069.09"
$ws.Range("E2").Value = "  -3.89%  "
$ws.Range("D3").Value = "1.596.33"
$ws.Range("E3").Value = "  -3.50%  "
$ws.Range("E4").Value = "  +0.51%  "
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "'301.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.50%  "
$ws.Range("D7").Value = "'0.3766"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.54%  "
$ws.Range("D8").Value = "'0.3637"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.11%  "
$ws.Range("D9").Value = "'47.67"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.99%  "
$ws.Range("D10").Value = "'1.002"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("D11").Value = "'1.271"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.44%  "
$ws.Range("D12").Value = "'0.08045"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.72%  "
$ws.Range("D13").Value = "'22.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.81%  "
$ws.Range("D14").Value = "'6.611"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.24%  "
$ws.Range("D15").Value = "'7.625"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.17%  "
$ws.Range("D16").Value = "'0.00001263"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.05%  "
$ws.Range("D17").Value = "1.596.35"
$ws.Range("E17").Value = "  -3.43%  "
$ws.Range("D18").Value = "'91.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.42%  "
$ws.Range("D19").Value = "'0.06784"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.70%  "
$ws.Range("D20").Value = "'18.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.58%  "
$ws.Range("D21").Value = "'6.577"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.60%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").Value = "'12.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.84%  "
$ws.Range("D24").Value = "23.110.61"
$ws.Range("E24").Value = "  -3.67%  "
$ws.Range("D25").Value = "'2.359"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.78%  "
$ws.Range("D26").Value = "'2.868"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.57%  "
$ws.Range("D27").Value = "'21.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.70%  "
$ws.Range("D28").Value = "'150.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.26%  "
$ws.Range("D29").Value = "'5.267"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.18%  "
$ws.Range("D30").Value = "'131.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.47%  "
$ws.Range("D31").Value = "'2.431"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("D32").Value = "'7.027"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.46%  "
$ws.Range("D33").Value = "1.770.29"
$ws.Range("E33").Value = "  -3.53%  "
$ws.Range("D34").Value = "'0.9859"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.88%  "
$ws.Range("D35").Value = "'0.07697"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.36%  "
$ws.Range("D36").Value = "'0.02778"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.53%  "
$ws.Range("D37").Value = "'6.269"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.83%  "
$ws.Range("D38").Value = "'0.2530"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.69%  "
$ws.Range("D39").Value = "'0.08859"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.04%  "
$ws.Range("D40").Value = "'10.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.79%  "
$ws.Range("E41").Value = "  -2.50%  "
$ws.Range("D42").Value = "'0.7147"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.57%  "
$ws.Range("D43").Value = "'12.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.85%  "
$ws.Range("D44").Value = "'15.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.20%  "
$ws.Range("D45").Value = "'0.6597"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.03%  "
$ws.Range("D46").Value = "'2.308"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.02%  "
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").Value = "'3.966"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.64%  "
$ws.Range("D51").Value = "'1.167"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.66%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.07976"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.09%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'131.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.33%  "
